# Automatische test-sync: 2025-06-22 18:45:50
#
# Adds a new incoming mail-log entry ("MVO-beleid") to the Logs sheet and
# refreshes the Dashboard's category-count table (which is kept sorted by
# descending count) to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new row (row 15)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(15, 1).Value = "MVO-beleid"
$logs.Cells.Item(15, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(15, 3).Value = "Hebben jullie een duurzaamheidsbeleid of MVO-doelen?"
$logs.Cells.Item(15, 4).Value = "Productinformatie"
$logs.Cells.Item(15, 6).Value = "2025-06-22 18:45:14"
$logs.Cells.Item(15, 7).Value = "Nee"

# Extend the conditional-formatting ranges on columns D and G so they keep
# covering the whole data range (D2:D14 -> D2:D15, G2:G14 -> G2:G15).
$dRange = $logs.Range("D2:D15")
$dConditions = $logs.Range("D2:D14").FormatConditions
for ($i = 1; $i -le $dConditions.Count; $i++) {
    $dConditions.Item($i).ModifyAppliesToRange($dRange)
}

$gRange = $logs.Range("G2:G15")
$gConditions = $logs.Range("G2:G14").FormatConditions
for ($i = 1; $i -le $gConditions.Count; $i++) {
    $gConditions.Item($i).ModifyAppliesToRange($gRange)
}

# ---------------------------------------------------------------------
# 2) Dashboard sheet: re-sort the category counts (descending by count)
#    and bump "Productinformatie" from 1 to 2.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Productinformatie"
$dash.Cells.Item(4, 2).Value = 2

$dash.Cells.Item(5, 1).Value = "Klacht / Probleem"
$dash.Cells.Item(5, 2).Value = 1

$dash.Cells.Item(6, 1).Value = "Openingstijden / Locatie"
$dash.Cells.Item(6, 2).Value = 1

$dash.Cells.Item(7, 1).Value = "Uitnodiging / Evenement"
$dash.Cells.Item(7, 2).Value = 1

$dash.Cells.Item(8, 1).Value = "Offerte / Prijsaanvraag"
$dash.Cells.Item(8, 2).Value = 1

$dash.Cells.Item(9, 1).Value = "Samenwerking / Partnerverzoek"
$dash.Cells.Item(9, 2).Value = 1

$dash.Cells.Item(10, 1).Value = "Afmelding / Nieuwsbrief"
$dash.Cells.Item(10, 2).Value = 1

$dash.Cells.Item(11, 1).Value = "Retour / Terugbetaling"
$dash.Cells.Item(11, 2).Value = 1
